$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 11:22"

# Refresh country ranking rows whose totals changed, causing reordering
$ws.Range("A13").Value = "Belgica"
$ws.Range("B13").Value = 23403
$ws.Range("C13").Value = 1209
$ws.Range("D13").Value = 4681
$ws.Range("E13").Value = 16482
$ws.Range("F13").Value = 1276
$ws.Range("G13").Value = 205
$ws.Range("H13").Value = 2240

$ws.Range("A14").Value = "Suiza"
$ws.Range("B14").Value = 22328
$ws.Range("C14").Value = 75
$ws.Range("D14").Value = 8704
$ws.Range("E14").Value = 12800
$ws.Range("F14").Value = 391
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 824

$ws.Range("A18").Value = "Austria"
$ws.Range("B18").Value = 12734
$ws.Range("C18").Value = 95
$ws.Range("D18").Value = 4512
$ws.Range("E18").Value = 7949
$ws.Range("F18").Value = 267
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 273

$ws.Range("A34").Value = "Malasia"
$ws.Range("B34").Value = 4119
$ws.Range("C34").Value = 156
$ws.Range("D34").Value = 1487
$ws.Range("E34").Value = 2567
$ws.Range("F34").Value = 76
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 65

$ws.Range("A35").Value = "Pakistan"
$ws.Range("B35").Value = 4072
$ws.Range("C35").Value = 37
$ws.Range("D35").Value = 467
$ws.Range("E35").Value = 3547
$ws.Range("F35").Value = 25
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 58

$ws.Range("A36").Value = "Ecuador"
$ws.Range("B36").Value = 3995
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 140
$ws.Range("E36").Value = 3635
$ws.Range("F36").Value = 156
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 220

$ws.Range("A54").Value = "Ucrania"
$ws.Range("B54").Value = 1668
$ws.Range("C54").Value = 206
$ws.Range("D54").Value = 35
$ws.Range("E54").Value = 1581
$ws.Range("F54").Value = 33
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 52

$ws.Range("A64").Value = "Eslovenia"
$ws.Range("B64").Value = 1091
$ws.Range("C64").Value = 32
$ws.Range("D64").Value = 120
$ws.Range("E64").Value = 931
$ws.Range("F64").Value = 35
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 40

$ws.Range("A65").Value = "Bielorrusia"
$ws.Range("B65").Value = 1066
$ws.Range("C65").Value = 205
$ws.Range("D65").Value = 77
$ws.Range("E65").Value = 976
$ws.Range("F65").Value = 33
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 13

$ws.Range("A71").Value = "Kuwait"
$ws.Range("B71").Value = 855
$ws.Range("C71").Value = 112
$ws.Range("D71").Value = 111
$ws.Range("E71").Value = 743
$ws.Range("F71").Value = 21
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 1

$ws.Range("A72").Value = "Barein"
$ws.Range("B72").Value = 811
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 465
$ws.Range("E72").Value = 341
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 5

$ws.Range("A73").Value = "Bosnia y Herzegovina"
$ws.Range("B73").Value = 777
$ws.Range("C73").Value = 13
$ws.Range("D73").Value = 77
$ws.Range("E73").Value = 667
$ws.Range("F73").Value = 4
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 33

$ws.Range("A119").Value = "Isla de Man"
$ws.Range("B119").Value = 158
$ws.Range("C119").Value = 8
$ws.Range("D119").Value = 80
$ws.Range("E119").Value = 77
$ws.Range("F119").Value = 7
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 1

$ws.Range("A120").Value = "Martinica"
$ws.Range("B120").Value = 152
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 50
$ws.Range("E120").Value = 98
$ws.Range("F120").Value = 20
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 4

$ws.Range("A146").Value = "Macao"
$ws.Range("B146").Value = 45
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 10
$ws.Range("E146").Value = 35
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0
